$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Writing a plain numeric-looking string (e.g. "541.15") via .Value lets Excel
    # auto-convert it to a Double, which loses the exact text (rounding / sci-notation).
    # Temporarily force a text format, assign, then clear the format again so the
    # cell's style index / formatting stays exactly as it was originally (unstyled),
    # while the underlying stored value remains the literal text from the source data.
    $range.NumberFormat = '@'
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range('D2').Value = '57.982.83'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '2.345.00'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').Value = '  +0.21%  '
Set-TextValue $ws.Range('D5') '541.15'
$ws.Range('E5').Value = '  -0.71%  '
Set-TextValue $ws.Range('D6') '134.12'
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('E7').Value = '  +0.23%  '
Set-TextValue $ws.Range('D8') '0.568'
$ws.Range('E8').Value = '  +5.84%  '
$ws.Range('E9').Value = '  +1.81%  '
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('E11').Value = '  -2.14%  '
Set-TextValue $ws.Range('D13') '23.73'
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('D14').Value = '2.761.26'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').Value = '57.923.89'
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').Value = '2.348.94'
$ws.Range('E17').Value = '  -0.20%  '
Set-TextValue $ws.Range('D18') '10.75'
$ws.Range('E18').Value = '  +1.58%  '
Set-TextValue $ws.Range('D19') '4.30'
$ws.Range('E19').Value = '  +2.08%  '
Set-TextValue $ws.Range('D20') '329.05'
$ws.Range('E20').Value = '  -1.78%  '
Set-TextValue $ws.Range('D21') '6.73'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('E23').Value = '  +1.87%  '
$ws.Range('E24').Value = '  -2.88%  '
Set-TextValue $ws.Range('D25') '0.998'
$ws.Range('E25').Value = '  +0.09%  '
Set-TextValue $ws.Range('D26') '8.32'
$ws.Range('E26').Value = '  -1.52%  '
$ws.Range('E27').Value = '  -6.18%  '
$ws.Range('E28').Value = '  -0.31%  '
Set-TextValue $ws.Range('D29') '170.46'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').Value = '0.0₃0734'
$ws.Range('E30').Value = '  -0.15%  '
Set-TextValue $ws.Range('D31') '6.13'
$ws.Range('E31').Value = '  -0.39%  '
Set-TextValue $ws.Range('D32') '18.33'
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('E33').Value = '  -1.95%  '
$ws.Range('E34').Value = '  +0.00%  '
Set-TextValue $ws.Range('D35') '0.999'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('E38').Value = '  -2.47%  '
Set-TextValue $ws.Range('D39') '39.11'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('E40').Value = '  -0.48%  '
Set-TextValue $ws.Range('D41') '289.94'
$ws.Range('E41').Value = '  +0.28%  '
Set-TextValue $ws.Range('D42') '140.46'
$ws.Range('E42').Value = '  -6.60%  '
Set-TextValue $ws.Range('D43') '3.64'
$ws.Range('E43').Value = '  +0.36%  '
Set-TextValue $ws.Range('D44') '0.0951'
$ws.Range('E44').Value = '  +2.36%  '
Set-TextValue $ws.Range('D45') '0.0511'
$ws.Range('E45').Value = '  +0.92%  '
Set-TextValue $ws.Range('D46') '18.93'
$ws.Range('E46').Value = '  -1.85%  '
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('E48').Value = '  +1.84%  '
Set-TextValue $ws.Range('D49') '0.381'
$ws.Range('E49').Value = '  -0.36%  '
Set-TextValue $ws.Range('D50') '11.08'
$ws.Range('E50').Value = '  +0.01%  '
Set-TextValue $ws.Range('D51') '4.71'
$ws.Range('E51').Value = '  +0.78%  '
